# Updates cryptos list per commit "Updated cryptos list on Wed Dec  6 08:15:16 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold plain decimal-number-looking text (e.g. "234.00") that must stay
# TEXT (matching the source inlineStr cells), not be auto-coerced to a Double by Excel and
# lose formatting (trailing zeros, etc). Prefixing with a literal leading apostrophe forces
# Excel to keep/interpret the entry as text, same as typing `234.00 in the Excel UI.

$ws.Range("D2").Value = "43.859.71"
$ws.Range("E2").Value = "  +5.27%  "
$ws.Range("D3").Value = "2.271.64"
$ws.Range("E3").Value = "  +2.97%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'234.00"
$ws.Range("E5").Value = "  +1.67%  "
$ws.Range("D6").Value = "'0.644"
$ws.Range("E6").Value = "  +4.63%  "
$ws.Range("D7").Value = "'64.04"
$ws.Range("E7").Value = "  +6.61%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "'0.431"
$ws.Range("E9").Value = "  +7.37%  "
$ws.Range("D10").Value = "'0.102"
$ws.Range("E10").Value = "  +15.31%  "
$ws.Range("D11").Value = "'57.47"
$ws.Range("E11").Value = "  +0.10%  "
$ws.Range("D12").Value = "'25.99"
$ws.Range("E12").Value = "  +16.92%  "
$ws.Range("D13").Value = "'0.104"
$ws.Range("E13").Value = "  +0.46%  "
$ws.Range("D14").Value = "2.607.92"
$ws.Range("E14").Value = "  +3.01%  "
$ws.Range("D15").Value = "'15.75"
$ws.Range("E15").Value = "  +2.64%  "
$ws.Range("D16").Value = "'5.95"
$ws.Range("E16").Value = "  +4.62%  "
$ws.Range("D17").Value = "'0.827"
$ws.Range("E17").Value = "  +4.37%  "
$ws.Range("D18").Value = "2.265.13"
$ws.Range("E18").Value = "  +2.29%  "
$ws.Range("D19").Value = "43.611.44"
$ws.Range("E19").Value = "  +4.95%  "
$ws.Range("D20").Value = "0.0₃0989"
$ws.Range("E20").Value = "  +9.89%  "
$ws.Range("D21").Value = "'74.16"
$ws.Range("E21").Value = "  +3.12%  "
$ws.Range("D22").Value = "'6.11"
$ws.Range("E22").Value = "  +0.80%  "
$ws.Range("D23").Value = "'249.52"
$ws.Range("E23").Value = "  +2.69%  "
$ws.Range("E24").Value = "  +0.14%  "
$ws.Range("D26").Value = "'2.32"
$ws.Range("E26").Value = "  +1.14%  "
$ws.Range("D27").Value = "'9.90"
$ws.Range("E27").Value = "  +2.51%  "
$ws.Range("D28").Value = "'173.15"
$ws.Range("E28").Value = "  +2.11%  "
$ws.Range("D29").Value = "'20.92"
$ws.Range("E29").Value = "  +6.01%  "
$ws.Range("E30").Value = "  -1.71%  "
$ws.Range("E31").Value = "  +1.46%  "
$ws.Range("E32").Value = "  +10.89%  "
$ws.Range("D33").Value = "'0.125"
$ws.Range("E33").Value = "  +2.88%  "
$ws.Range("D34").Value = "'0.0688"
$ws.Range("E34").Value = "  +6.32%  "
$ws.Range("E35").Value = "  +1.96%  "
$ws.Range("D36").Value = "'4.76"
$ws.Range("E36").Value = "  +2.93%  "
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").Value = "'3.88"
$ws.Range("E37").Value = "  +9.78%  "
$ws.Range("B38").Value = "THORChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D38").Value = "'6.86"
$ws.Range("E38").Value = "  +6.41%  "
$ws.Range("D39").Value = "'2.35"
$ws.Range("E39").Value = "  -0.82%  "
$ws.Range("D40").Value = "'0.0250"
$ws.Range("E40").Value = "  +5.50%  "
$ws.Range("E41").Value = "  +0.16%  "
$ws.Range("D42").Value = "'8.41"
$ws.Range("E42").Value = "  -1.65%  "
$ws.Range("D43").Value = "'17.53"
$ws.Range("E43").Value = "  +6.95%  "
$ws.Range("D44").Value = "'10.50"
$ws.Range("E44").Value = "  +21.90%  "
$ws.Range("D45").Value = "'0.0966"
$ws.Range("E45").Value = "  +0.97%  "
$ws.Range("D46").Value = "'4.47"
$ws.Range("E46").Value = "  +2.85%  "
$ws.Range("E47").Value = "  +0.50%  "
$ws.Range("D48").Value = "'97.86"
$ws.Range("E48").Value = "  +0.60%  "
$ws.Range("D49").Value = "1.481.87"
$ws.Range("E49").Value = "  +0.97%  "
$ws.Range("D50").Value = "'2.33"
$ws.Range("E50").Value = "  +4.18%  "
$ws.Range("E51").Value = "  +0.99%  "
